$d = $word.ActiveDocument

# --- 1. Fill in the third data row of the Resource Requirements Table ---
# (table 1 in the document: Item ID / Description / Unit Cost / Quantity /
#  Subtotal / Event ID / Dependencies). Row 4 is currently all-empty cells;
#  populate the first five columns, leaving Event ID / Dependencies blank.
$reqTable = $d.Tables.Item(1)
$reqTable.Cell(4, 1).Range.Text = "3"
$reqTable.Cell(4, 2).Range.Text = "Laptops"
$reqTable.Cell(4, 3).Range.Text = "`$1,000"
$reqTable.Cell(4, 4).Range.Text = "3"
$reqTable.Cell(4, 5).Range.Text = "`$3,000.00"

# --- 2. Append two new paragraphs at the very end of the document body,
#         right before the trailing paragraph mark that precedes sectPr ---

# 2a. A blank paragraph.
$insertAt = $d.Content.End - 1
$rng = $d.Range($insertAt, $insertAt)
$rng.Text = "`r"

# 2b. A paragraph with the cost-justification text.
$insertAt = $d.Content.End - 1
$rng = $d.Range($insertAt, $insertAt)
$rng.Text = "`r"

$insertAt = $d.Content.End - 1
$rng = $d.Range($insertAt, $insertAt)
$rng.InsertAfter("Laptops must be upgraded every two years, but must be purchased immediately at `$1,000 per laptop for all three team members.")

$insertAt = $d.Content.End - 1
$rng = $d.Range($insertAt, $insertAt)
$rng.InsertAfter("  Sprites must also be purchased from a graphic artist for all units and items in the game.  `$10.00 per sprite is around the price of a custom sprite designed by an artist that is not available for use in existing or future games.")

$insertAt = $d.Content.End - 1
$rng = $d.Range($insertAt, $insertAt)
$rng.InsertAfter("  The hourly rate of `$30.00 for each programmer and `$50.00 for the team leader is a market average for medium-level programmers and lead programmers on an i")

$insertAt = $d.Content.End - 1
$rng = $d.Range($insertAt, $insertAt)
$rng.InsertAfter("ndependent game.  The total cost for the project will be: ")
